$wb = $excel.ActiveWorkbook

# Clear the stray empty cells in the "ODI Batting" sheet
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("B2").ClearContents()
$odiBatting.Range("B3").ClearContents()

# Add the new "ODI Batting Extra" sheet at the end
$extra = $wb.Worksheets.Add()
$extra.Name = "ODI Batting Extra"

$extra.Cells.Item(1, 1).Value = "MATCH_CODE"
$extra.Cells.Item(1, 2).Value = "BATTING_POSITION"
$extra.Cells.Item(1, 3).Value = "NUM_4"
$extra.Cells.Item(1, 4).Value = "NUM_6"
$extra.Cells.Item(1, 5).Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Cells.Item(1, 6).Value = "MAN_OF_MATCH"

$extra.Cells.Item(2, 1).Value = "4237"
$extra.Cells.Item(2, 2).Value = 10
$extra.Cells.Item(2, 6).Value = "NO"

$extra.Cells.Item(3, 1).Value = "4238"
